$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1101.6666
$ws.Range("I41").Value = 1253
$ws.Range("J41").Value = 572
$ws.Range("K41").Value = 1253
$ws.Range("L41").Value = 572
$ws.Range("M41").Value = -813
$ws.Range("N41").Value = -1452

$ws.Range("H43").Value = 4625.1064
$ws.Range("J43").Value = 4250.364
$ws.Range("L43").Value = 4250.364
$ws.Range("N43").Value = -4388.364

$ws.Range("H116").Value = 55955.55
$ws.Range("I116").Value = 146763.86
$ws.Range("J116").Value = 7058.769
$ws.Range("K116").Value = 146763.86
$ws.Range("L116").Value = 7058.769
$ws.Range("M116").Value = -143321.86
$ws.Range("N116").Value = -13942.769

$ws.Range("H132").Value = 33329.605
$ws.Range("I132").Value = 37493.516
$ws.Range("J132").Value = 3141.25
$ws.Range("K132").Value = 112480.548
$ws.Range("L132").Value = 9423.75
$ws.Range("M132").Value = -109950.548
$ws.Range("N132").Value = -14483.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 623
$ws.Range("J17").Value = 354.5
$ws.Range("L17").Value = 354.5
$ws.Range("N17").Value = -700.5

$ws.Range("H110").Value = 4515.4585
$ws.Range("I110").Value = 3972.8667
$ws.Range("J110").Value = 5419.778
$ws.Range("K110").Value = 3972.8667
$ws.Range("L110").Value = 5419.778
$ws.Range("M110").Value = -1927.8667
$ws.Range("N110").Value = -9509.778

$ws.Range("H132").Value = 31257108
$ws.Range("I132").Value = 4761.4
$ws.Range("J132").Value = 83344350
$ws.Range("K132").Value = 14284.2
$ws.Range("L132").Value = 250033050
$ws.Range("M132").Value = -11754.2
$ws.Range("N132").Value = -250038110

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1893.6818
$ws.Range("I58").Value = 1874.45
$ws.Range("J58").Value = 2086
$ws.Range("K58").Value = 1874.45
$ws.Range("L58").Value = 2086
$ws.Range("M58").Value = -1671.45
$ws.Range("N58").Value = -2492

$ws.Range("H124").Value = 88331.664
$ws.Range("J124").Value = 88331.664
$ws.Range("L124").Value = 88331.664
$ws.Range("N124").Value = -93241.664

$ws.Range("H136").Value = 1893.6818
$ws.Range("I136").Value = 1874.45
$ws.Range("J136").Value = 2086
$ws.Range("K136").Value = 5623.35
$ws.Range("L136").Value = 6258
$ws.Range("M136").Value = -3073.35
$ws.Range("N136").Value = -11358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1331

$ws.Range("H27").Value = 500
$ws.Range("I27").Value = 500
$ws.Range("K27").Value = 1500
$ws.Range("M27").Value = -1398

$ws.Range("H29").Value = 1938.25
$ws.Range("I29").Value = 499
$ws.Range("K29").Value = 1497
$ws.Range("M29").Value = -1220

$ws.Range("H33").Value = 240.63637
$ws.Range("I33").Value = 90.5
$ws.Range("J33").Value = 420.8
$ws.Range("K33").Value = 543
$ws.Range("L33").Value = 2524.8
$ws.Range("M33").Value = -260
$ws.Range("N33").Value = -3090.8

$ws.Range("H36").Value = 1209
$ws.Range("J36").Value = 1003
$ws.Range("L36").Value = 3009
$ws.Range("N36").Value = -3347

$ws.Range("H38").Value = 103.35
$ws.Range("J38").Value = 168.27272
$ws.Range("L38").Value = 504.81816
$ws.Range("N38").Value = -1198.81816

$ws.Range("H44").Value = 823.8333
$ws.Range("I44").Value = 688.2
$ws.Range("J44").Value = 1502
$ws.Range("K44").Value = 2064.6
$ws.Range("L44").Value = 4506
$ws.Range("M44").Value = -1666.6
$ws.Range("N44").Value = -5302

$ws.Range("H48").Value = 11722
$ws.Range("J48").Value = 11722
$ws.Range("L48").Value = 35166
$ws.Range("N48").Value = -35666

$ws.Range("H63").Value = 2222
$ws.Range("J63").Value = 2222
$ws.Range("L63").Value = 6666
$ws.Range("N63").Value = -8164

$ws.Range("H66").Value = 2222
$ws.Range("J66").Value = 2222
$ws.Range("L66").Value = 19998
$ws.Range("N66").Value = -27486

$ws.Range("H109").Value = 3432.9167
$ws.Range("I109").Value = 2619.5
$ws.Range("K109").Value = 7858.5
$ws.Range("M109").Value = -6818.5

$ws.Range("H132").Value = 2337.6
$ws.Range("I132").Value = 1895
$ws.Range("K132").Value = 17055
$ws.Range("M132").Value = -14525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H46").Value = 16020.5
$ws.Range("I46").Value = 16020.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 16020.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -15864.5
$ws.Range("N46").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H82").Value = 100000
$ws.Range("J82").Value = 100000
$ws.Range("L82").Value = 100000
$ws.Range("N82").Value = -100766

$ws.Range("H85").Value = 100000
$ws.Range("J85").Value = 100000
$ws.Range("L85").Value = 100000
$ws.Range("N85").Value = -102652

$ws.Range("H132").Value = 2345.2144
$ws.Range("I132").Value = 2345.2144
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7035.6432
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4505.6432
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888

$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
